$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CPPbI")

# New header cells for column B and C
$ws.Range("B1").Value = "energy related emissions"
$ws.Range("C1").Value = "process emissions"

# Column widths
$ws.Columns.Item(1).ColumnWidth = 43.7109375
$ws.Columns.Item(2).ColumnWidth = 24.7109375
$ws.Columns.Item(3).ColumnWidth = 25.85546875

# Populate column C with same values as column B (rows 2-9)
$ws.Range("C2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("C4").Value = 1
$ws.Range("C5").Value = 1
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 1
